$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-224) holds the "Förändrad" date serial value.
# It needs to move from 45171 (2023-09-02) to 45172 (2023-09-03) for every row.
$ws.Range("C2:C224").Value = 45172
